$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two users who moved away (ROUEN agency staff no longer
# distributed in this update): row 4 = FOUCAULT, row 6 = JEANNE.
# Delete row 6 first so row 4's index is not affected by the shift.
$ws.Rows(6).Delete() | Out-Null
$ws.Rows(4).Delete() | Out-Null

# Renumber the ID column (A) sequentially for the remaining records.
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7

# The hyperlinked e-mail (FRESNAIS Mathis) used to live on row 5 and is
# now on row 4; move the hyperlink to follow the cell.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("Q4"), "mailto:m.fresnais@ouestisol.fr") | Out-Null
$ws.Range("Q4").Style = "Lien hypertexte"

# Restore the selection to match the refreshed data range.
$ws.Range("A2:A8").Select() | Out-Null
